# PivotTableShowDataAs.xlsx — "Add support for '% of' show data as setting option"
#
# The refreshed pivot table at B30:J41 ("PivotTable3") now shows rolled-up
# values on its (previously blank) group-header rows for San Francisco (33),
# Chicago (35) and Nashville (38) - the normal result of Excel recomputing
# the "% of" / count aggregates for those outline rows on refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PivotTables")

# Row 33 - San Francisco header row
$ws.Range("C33").Value = 0.94329990470572223
$ws.Range("D33").Value = 0.40344492964580303
$ws.Range("E33").Value = 0.80767362797474507
$ws.Range("F33").Value = 2
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 1
$ws.Range("I33").Value = 0.62802431029360672
$ws.Range("J33").Value = 5

# Row 35 - Chicago header row (F35 stays blank)
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 0.193110140708394
$ws.Range("E35").Value = 0.19232637202525499
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 1
$ws.Range("I35").Value = 0.15005110801162141
$ws.Range("J35").Value = 7

# Row 38 - Nashville header row (H38 stays blank)
$ws.Range("C38").Value = 5.6700095294277802E-2
$ws.Range("D38").Value = 0.40344492964580303
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 2
$ws.Range("I38").Value = 0.22192458169477189
$ws.Range("J38").Value = 3
